$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Extend header row with two more columns (P1, Q1), copying O1's style (bold, centered, bordered)
$ws.Range("O1").Copy($ws.Range("P1"))
$ws.Range("O1").Copy($ws.Range("Q1"))
$ws.Range("P1").Value = 14
$ws.Range("Q1").Value = 15

# Update data rows 2..25: columns I, K, M, O change values, and new columns P, Q are appended
for ($r = 2; $r -le 25; $r++) {
    $ws.Cells.Item($r, 9).Value = 2   # I -> 2
    $ws.Cells.Item($r, 11).Value = 1  # K -> 1
    $ws.Cells.Item($r, 13).Value = 2  # M -> 2
    $ws.Cells.Item($r, 15).Value = 1  # O -> 1
    $ws.Cells.Item($r, 16).Value = 2  # P -> 2
    $ws.Cells.Item($r, 17).Value = 2  # Q -> 2
}
